$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "Item"
$ws.Cells.Item(1, 2).Value = "Descrição"
$ws.Cells.Item(1, 3).Value = "Marca"
$ws.Cells.Item(1, 4).Value = "Registro"
$ws.Cells.Item(1, 5).Value = "Apresentação"
$ws.Cells.Item(1, 6).Value = "PDF"

# Data rows: Item, Descrição, Marca, Registro, Apresentação, PDF
$data = @(
    @(1,  "HEPARINA SODICA SUBCUT 5000UI", "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980371", "5000 UI/ML SOL INJ CX 25 FA VD INC X 5 ML ATIVA", "OK"),
    @(4,  "PROMETAZINA 25MG, CLORIDRATO", "LABORATÓRIO TEUTO BRASILEIRO S/A", "103700321", "25 MG COM REV BL AL PLAS INC X 200 (EMB. HOSP.) ATIVA", "OK"),
    @(5,  "HALOPERIDOL 5MG", "CELLERA FARMACÊUTICA S.A.", "112360011", "5 MG COM CT BL AL PLAS TRANS X 20 ATIVA", "OK"),
    @(6,  "CLORPROMAZINA 40MG/ML SOL ORAL", "SANOFI MEDLEY FARMACÊUTICA LTDA", "Último registro encontrado: 183260385", "Não encontrado", "Pendente"),
    @(7,  "HALOPERIDOL 2MG/ML SOL ORAL", "CELLERA FARMACÊUTICA S.A.", "112360011", "2 MG/ML SOL GOT OR CT FR GOT PLAS OPC X 30 ML ATIVA", "OK"),
    @(8,  "CLORPROMAZINA 25MG", "SANOFI MEDLEY FARMACÊUTICA LTDA", "183260385", "25 MG COM REV CT BL AL PLAS OPC X 20 ATIVA", "OK"),
    @(9,  "CODEINA 30MG", "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980199", "30 MG COM CT BL AL PLAS TRANS X 30 ATIVA", "OK"),
    @(10, "IMIPRAMINA 25MG", "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980023", "25 MG COM REV CX 20 ENV AL POLIET X 10 (EMB HOSP) CANCELADA OU CADUCA", "OK"),
    @(11, "RISPERIDONA 3MG", "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA", "103920197", "3 MG COM REV CT BL AL PLAS TRANS X 10 ATIVA", "OK"),
    @(13, "RISPERIDONA 1MG", "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA", "103920197", "1 MG COM REV CT BL AL PLAS TRANS X 10 ATIVA", "OK"),
    @(14, "LEVOMEPROMAZINA 4% GOTAS", "SANOFI MEDLEY FARMACÊUTICA LTDA", "Último registro encontrado: 183260316", "Não encontrado", "Pendente"),
    @(16, "LIDOCAINA 2% C/ VASO CONSTRITO", "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980249", "2% GEL TOP CT BG AL X 30 G + APLIC ATIVA", "OK"),
    @(17, "NITRATO DE CERIO +SULFADIAZINA", "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "Último registro encontrado: 102980560", "Não encontrado", "Pendente"),
    @(18, "COLAGENASE+CLORAFENICOL POMADA 30g", "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA", "102980431", "0,6 U/G POM DERM CT 01 BG AL X 30 G + ESP PLAS ATIVA", "OK")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    # Registro column holds numeric-looking codes that must stay as text
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}
